$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.016.41"
$ws.Range("E2").Value = "  -1.14%  "

$ws.Range("D3").Value = "2.353.97"
$ws.Range("E3").Value = "  -1.40%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.59"
$ws.Range("E5").Value = "  +0.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.46"
$ws.Range("E6").Value = "  -1.57%  "

$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("E8").Value = "  -2.43%  "

$ws.Range("D9").Value = "2.367.30"
$ws.Range("E9").Value = "  -1.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0974"
$ws.Range("E10").Value = "  +0.11%  "

$ws.Range("E11").Value = "  -0.57%  "

$ws.Range("E12").Value = "  +2.87%  "

$ws.Range("E13").Value = "  -1.17%  "

$ws.Range("D14").Value = "2.771.60"
$ws.Range("E14").Value = "  -1.42%  "

$ws.Range("D15").Value = "55.960.14"
$ws.Range("E15").Value = "  -1.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.53"
$ws.Range("E16").Value = "  -0.54%  "

$ws.Range("E17").Value = "  -0.38%  "

$ws.Range("D18").Value = "2.398.57"
$ws.Range("E18").Value = "  +0.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.93"
$ws.Range("E19").Value = "  -2.35%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "311.09"
$ws.Range("E20").Value = "  +0.76%  "

$ws.Range("E21").Value = "  -0.67%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.19"
$ws.Range("E22").Value = "  -1.01%  "

$ws.Range("E23").Value = "  -0.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.42"
$ws.Range("E24").Value = "  -0.87%  "

$ws.Range("E25").Value = "  -0.44%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.370"
$ws.Range("E26").Value = "  -2.96%  "

$ws.Range("E27").Value = "  -2.69%  "

$ws.Range("E28").Value = "  -2.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.62"
$ws.Range("E29").Value = "  -2.55%  "

$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0708"
$ws.Range("E30").Value = "  -2.39%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.64"
$ws.Range("E31").Value = "  -0.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.75"
$ws.Range("E33").Value = "  -1.79%  "

$ws.Range("E34").Value = "  -0.20%  "

$ws.Range("E35").Value = "  -4.27%  "

$ws.Range("E36").Value = "  -0.52%  "

$ws.Range("E37").Value = "  -1.51%  "

$ws.Range("E38").Value = "  +2.69%  "

$ws.Range("E39").Value = "  -4.18%  "

$ws.Range("E40").Value = "  -1.60%  "

$ws.Range("E41").Value = "  -2.99%  "

$ws.Range("E42").Value = "  -0.72%  "

$ws.Range("E43").Value = "  +0.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "125.40"
$ws.Range("E44").Value = "  -4.09%  "

$ws.Range("E45").Value = "  -1.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "242.77"
$ws.Range("E47").Value = "  -1.89%  "

$ws.Range("E48").Value = "  -0.87%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.79"
$ws.Range("E49").Value = "  -1.12%  "

$ws.Range("E50").Value = "  -1.50%  "

$ws.Range("E51").Value = "  -3.06%  "

